# NOTE: the order in which shared strings get touched below matters - it
# reproduces the exact shared-string table ordering of the target file
# (old "company adoption fraction" slot dropped, old last slot "operation
# cost" renamed in place to "minimum operation cost", then two brand new
# strings appended: "initial push time", "initial company adoption
# fraction").
$wb = $excel.ActiveWorkbook

# --- "cost" sheet --------------------------------------------------------
$ws4 = $wb.Worksheets.Item("cost")

# Rename the "operation cost" label to "minimum operation cost" (value
# unchanged, 1000).
$ws4.Range("A4").Value = "minimum operation cost"

# --- "company" sheet ---------------------------------------------------
$ws1 = $wb.Worksheets.Item("company")

# Add new row 4 ("initial push time").
$ws1.Range("A4").Value = "initial push time"
$ws1.Range("B4").Value = 36

# Rename the row-3 label (value itself, 0.3, is unchanged).
$ws1.Range("A3").Value = "initial company adoption fraction"

# --- Selection / active-cell bookkeeping (cosmetic, mirrors the diff) ---
$ws1.Range("B12").Select()

$ws3 = $wb.Worksheets.Item("demand")
$ws3.Range("B3").Select()

$ws4.Range("B1").Select()

# --- Window size bookkeeping (cosmetic; windowHeight 16935 -> 9420) -----
try {
    $win = $wb.Windows.Item(1)
    $win.Height = 9420
} catch {
    # Not fatal if unavailable in this host - purely cosmetic window chrome.
}
